$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H, matching the formatting of the
# existing header cells (e.g. G1 - bold, bordered, centered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the Save values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
